# spring 23 week 7 inputs
# Append 20 new data rows (776-795) to the "Nine" sheet, matching
# the pattern of the existing A:D matchup data columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

$newRows = @(
    @(5,16,4,4),
    @(5,13,4,7),
    @(2,17,5,3),
    @(3,14,4,6),
    @(6,17,5,3),
    @(5,13,3,7),
    @(5,12,4,8),
    @(3,17,2,3),
    @(5,12,7,8),
    @(6,13,5,7),
    @(4,12,5,8),
    @(6,8,7,12),
    @(5,7,4,13),
    @(8,15,4,5),
    @(1,14,4,6),
    @(3,8,5,12),
    @(5,8,4,12),
    @(7,8,4,12),
    @(3,17,4,3),
    @(5,4,4,16)
)

$startRow = 776

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

$lastRow = $startRow + $newRows.Count - 1
$nextRow = $lastRow + 1

$ws.Range("A$nextRow").Select()
